$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows of data right after the header row (old row2 -> new row4, etc.)
# and append eight new rows at the bottom (new rows 24-31).
# Since the runtime easily supports direct cell writes, rewrite rows 2-31 in place
# to their final values (equivalent to inserting 2 rows at the top and 8 at the bottom).

$ws.Cells.Item(2,1).Value = -0.3733085989952087
$ws.Cells.Item(2,2).Value = -2.073603868484497
$ws.Cells.Item(2,3).Value = 2.090381860733032
$ws.Cells.Item(3,1).Value = 0.200868934392929
$ws.Cells.Item(3,2).Value = -1.047221541404724
$ws.Cells.Item(3,3).Value = 2.184124946594238
$ws.Cells.Item(4,1).Value = -1.194094896316528
$ws.Cells.Item(4,2).Value = 0.9634650945663452
$ws.Cells.Item(4,3).Value = 0.0717056170105934
$ws.Cells.Item(5,1).Value = -1.414604663848877
$ws.Cells.Item(5,2).Value = 0.444681316614151
$ws.Cells.Item(5,3).Value = 1.771335005760193
$ws.Cells.Item(6,1).Value = -1.245760202407837
$ws.Cells.Item(6,2).Value = 3.215966701507568
$ws.Cells.Item(6,3).Value = 5.047555923461914
$ws.Cells.Item(7,1).Value = -0.3136537969112396
$ws.Cells.Item(7,2).Value = 3.28307843208313
$ws.Cells.Item(7,3).Value = 0.4722450375556946
$ws.Cells.Item(8,1).Value = 2.299439907073975
$ws.Cells.Item(8,2).Value = 0.4132560193538666
$ws.Cells.Item(8,3).Value = -0.0342882014811039
$ws.Cells.Item(9,1).Value = 2.066679716110229
$ws.Cells.Item(9,2).Value = -1.99637222290039
$ws.Cells.Item(9,3).Value = -4.815195083618164
$ws.Cells.Item(10,1).Value = 3.397727489471436
$ws.Cells.Item(10,2).Value = 7.844008445739746
$ws.Cells.Item(10,3).Value = -6.625612258911133
$ws.Cells.Item(11,1).Value = -3.216499328613281
$ws.Cells.Item(11,2).Value = -3.29652738571167
$ws.Cells.Item(11,3).Value = 5.164734840393066
$ws.Cells.Item(12,1).Value = -6.98007869720459
$ws.Cells.Item(12,2).Value = 0.7972838878631592
$ws.Cells.Item(12,3).Value = 8.478240013122559
$ws.Cells.Item(13,1).Value = -2.500109195709229
$ws.Cells.Item(13,2).Value = 7.255982398986816
$ws.Cells.Item(13,3).Value = 6.759435653686523
$ws.Cells.Item(14,1).Value = 4.161522388458252
$ws.Cells.Item(14,2).Value = -0.7580022215843201
$ws.Cells.Item(14,3).Value = -2.298507690429688
$ws.Cells.Item(15,1).Value = 3.58042049407959
$ws.Cells.Item(15,2).Value = -4.013983249664307
$ws.Cells.Item(15,3).Value = -4.726245403289795
$ws.Cells.Item(16,1).Value = 0.5662546157836914
$ws.Cells.Item(16,2).Value = 0.348807543516159
$ws.Cells.Item(16,3).Value = -7.386210918426514
$ws.Cells.Item(17,1).Value = -5.530786514282227
$ws.Cells.Item(17,2).Value = -3.427022218704224
$ws.Cells.Item(17,3).Value = 3.131677865982056
$ws.Cells.Item(18,1).Value = -6.868226051330566
$ws.Cells.Item(18,2).Value = -1.806755065917969
$ws.Cells.Item(18,3).Value = 10.58639812469482
$ws.Cells.Item(19,1).Value = 1.907955169677734
$ws.Cells.Item(19,2).Value = 14.34744739532471
$ws.Cells.Item(19,3).Value = 0.4956808686256408
$ws.Cells.Item(20,1).Value = -0.698746919631958
$ws.Cells.Item(20,2).Value = -6.705373764038086
$ws.Cells.Item(20,3).Value = 1.454951524734497
$ws.Cells.Item(21,1).Value = 6.464091300964356
$ws.Cells.Item(21,2).Value = -1.321527123451233
$ws.Cells.Item(21,3).Value = -5.915080547332764
$ws.Cells.Item(22,1).Value = 4.413990020751953
$ws.Cells.Item(22,2).Value = -3.110638856887817
$ws.Cells.Item(22,3).Value = -1.591172218322754
$ws.Cells.Item(23,1).Value = -1.979194760322571
$ws.Cells.Item(23,2).Value = 1.858819842338562
$ws.Cells.Item(23,3).Value = 3.587078332901001
$ws.Cells.Item(24,1).Value = -9.171860694885254
$ws.Cells.Item(24,2).Value = -10.49505233764648
$ws.Cells.Item(24,3).Value = -3.089466571807861
$ws.Cells.Item(25,1).Value = 5.832921981811523
$ws.Cells.Item(25,2).Value = -5.742907047271729
$ws.Cells.Item(25,3).Value = 5.470998287200928
$ws.Cells.Item(26,1).Value = 0.7526758909225464
$ws.Cells.Item(26,2).Value = -7.854794025421143
$ws.Cells.Item(26,3).Value = 0.5755757093429565
$ws.Cells.Item(27,1).Value = 8.307531356811523
$ws.Cells.Item(27,2).Value = 6.279134750366211
$ws.Cells.Item(27,3).Value = -0.9227187633514404
$ws.Cells.Item(28,1).Value = -3.3331458568573
$ws.Cells.Item(28,2).Value = 3.669769525527954
$ws.Cells.Item(28,3).Value = -1.161337971687317
$ws.Cells.Item(29,1).Value = -1.678790211677551
$ws.Cells.Item(29,2).Value = 2.915562152862549
$ws.Cells.Item(29,3).Value = 3.480551958084106
$ws.Cells.Item(30,1).Value = -8.003265380859375
$ws.Cells.Item(30,2).Value = 1.636179566383362
$ws.Cells.Item(30,3).Value = -1.630586981773376
$ws.Cells.Item(31,1).Value = 1.040297269821167
$ws.Cells.Item(31,2).Value = -2.371345281600952
$ws.Cells.Item(31,3).Value = 6.342917442321777